$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.311.09"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.875.12"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - XRP
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7108"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.09%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.08%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07859"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.05%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3129"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.74%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.74%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08387"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.03%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.872.21"
$ws.Range("E12").Value = "  -2.29%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.244"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.61%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7175"
$ws.Range("D14").ClearFormats()

# Row 15 - Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.23"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.00%  "

# Row 16 - Uniswap
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.212"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.10%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008341"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.02%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.314.19"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.72%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.33%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.124.70"
$ws.Range("E21").Value = "  -1.36%  "

# Row 22 - Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.13%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.786"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.59%  "

# Row 24 - BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.04%  "

# Row 25 - Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.53%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.054"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.48%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.48%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.15%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.24%  "

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.422"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.12%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.349"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.77%  "

# Row 32 - Toncoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.210"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.07%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05357"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.16%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.947"
$ws.Range("D34").ClearFormats()

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.53%  "

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7496"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.37%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.693"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.43%  "

# Row 38 - Maker
$ws.Range("D38").Value = "1.297.59"
$ws.Range("E38").Value = "  +12.20%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +1.59%  "

# Row 40 - MXToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.736"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.76%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.570"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.17%  "

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.64%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8957"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.16%  "

# Row 44 - Aave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - BabyDogeCoin
$ws.Range("E45").Value = "  +8.91%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  +0.06%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "2.021.23"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48 - RenderToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.802"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.00%  "

# Row 49 - Mantle
$ws.Range("E49").Value = "  +0.11%  "

# Row 50 - EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.457"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.79%  "

# Row 51 - TheSandbox
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4359"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.47%  "
